$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155; existing rows 155..182 shift down to 156..183,
# preserving the row-155 date-column style (s="2") that Excel carries into the
# newly-opened row.
$ws.Rows.Item(155).Insert()

# Populate the new row 155 with the new data record.
$ws.Cells.Item(155, 1).Value  = 3
$ws.Cells.Item(155, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value  = "Coquimbo"
$ws.Cells.Item(155, 4).Value  = 44522
$ws.Cells.Item(155, 5).Value  = 5
$ws.Cells.Item(155, 6).Value  = 100112001
$ws.Cells.Item(155, 7).Value  = "Berenjena"
$ws.Cells.Item(155, 8).Value  = "Sin especificar"
$ws.Cells.Item(155, 9).Value  = "Primera"
$ws.Cells.Item(155, 10).Value = 140
$ws.Cells.Item(155, 11).Value = 7000
$ws.Cells.Item(155, 12).Value = 7500
$ws.Cells.Item(155, 13).Value = 7250
$ws.Cells.Item(155, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 121
$ws.Cells.Item(155, 17).Value = 60
$ws.Cells.Item(155, 18).Value = "Hortaliza"
